$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Text Modify")

# Move B4 (PROPER formula) -> B5
$ws.Range("B5").Formula = $ws.Range("B4").Formula
$ws.Range("B4").ClearContents()

# Move B6 (LOWER formula) -> B8
$ws.Range("B8").Formula = $ws.Range("B6").Formula
$ws.Range("B6").ClearContents()

# New row 3: A3 blank but styled like A1/A2; B3 shows the literal formula text "=UPPER(A1)"
$ws.Range("A1").Copy()
$ws.Range("A3").PasteSpecial(-4122)
$ws.Range("B3").Value = "'=UPPER(A1)"
$ws.Range("B3").Borders.Weight = -4138

# New row 6: B6 shows the literal formula text "=PROPER(A1)"
$ws.Range("B6").Value = "'=PROPER(A1)"
$ws.Range("B6").Borders.Weight = -4138

# New row 9: B9 shows the literal formula text "=LOWER(A1)"
$ws.Range("B9").Value = "'=LOWER(A1)"
$ws.Range("B9").Borders.Weight = -4138

# Row heights
$ws.Rows.Item(5).RowHeight = 17
$ws.Rows.Item(6).RowHeight = 17
$ws.Rows.Item(8).RowHeight = 17
$ws.Rows.Item(9).RowHeight = 17

# Column B widened to fit the new longer text
$ws.Columns.Item(2).ColumnWidth = 57

# Selection moves to B14
$ws.Range("B14").Select()

Write-Output "done"
